$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.893.47"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.872.80"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.033"
$ws.Range("E4").Value = "  +0.68%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.77"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.030"
$ws.Range("E6").Value = "  +0.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4427"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3833"
$ws.Range("E8").Value = "  +3.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07471"
$ws.Range("E9").Value = "  +1.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8910"
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("D12").Value = "1.884.58"
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.591"
$ws.Range("E13").Value = "  +2.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.789"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07201"
$ws.Range("E15").Value = "  +1.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "86.14"
$ws.Range("E16").Value = "  +4.79%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.036"
$ws.Range("E17").Value = "  +0.74%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009154"
$ws.Range("E18").Value = "  +1.86%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.029"
$ws.Range("E19").Value = "  +0.60%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.64"
$ws.Range("E20").Value = "  +2.19%  "
$ws.Range("D21").Value = "27.916.14"
$ws.Range("E21").Value = "  +1.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.337"
$ws.Range("E22").Value = "  +1.90%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.33"
$ws.Range("E23").Value = "  +1.74%  "
$ws.Range("D24").Value = "2.098.51"
$ws.Range("E24").Value = "  -1.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.026"
$ws.Range("E25").Value = "  +6.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "158.77"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.93"
$ws.Range("E27").Value = "  +2.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.427"
$ws.Range("E28").Value = "  +3.52%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.993"
$ws.Range("E29").Value = "  +4.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.48"
$ws.Range("E30").Value = "  +2.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09071"
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.233"
$ws.Range("E32").Value = "  +2.91%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7857"
$ws.Range("E33").Value = "  +3.57%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.025"
$ws.Range("E34").Value = "  +5.75%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.617"
$ws.Range("E35").Value = "  +3.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.031"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.152"
$ws.Range("E37").Value = "  +0.88%  "
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01994"
$ws.Range("E38").Value = "  +1.84%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05374"
$ws.Range("E39").Value = "  +2.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.885"
$ws.Range("E40").Value = "  +3.54%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5251"
$ws.Range("E41").Value = "  +1.98%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1700"
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.961"
$ws.Range("E43").Value = "  +6.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.899"
$ws.Range("E44").Value = "  +4.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "112.45"
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.78"
$ws.Range("E46").Value = "  +2.16%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.06630"
$ws.Range("E47").Value = "  +5.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.728"
$ws.Range("E48").Value = "  +3.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.032"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4764"
$ws.Range("E50").Value = "  +3.17%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.941"
$ws.Range("E51").Value = "  +3.45%  "
